$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C (product image URL column) to match the new, longer hyperlink text.
$ws.Columns.Item(3).ColumnWidth = 89.83

# Every value in column D (rows 1-70) is a relative "/static/img/..." path.
# Turn each one into an absolute https://scrapingclub.com/... hyperlink, and
# make the displayed cell text the full absolute URL as well.
for ($r = 1; $r -le 70; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $old = $cell.Value2
    $url = "https://scrapingclub.com" + $old
    $ws.Hyperlinks.Add($cell, $url, [Type]::Missing, [Type]::Missing, $url) | Out-Null
}

Write-Output "done"
